# Add a "double header" row above the existing header row:
# an extra row that carries a couple of category-style labels
# ("Extra header" over column A, "Input category" over column I),
# pushing the old header + data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; existing rows 1-2 become 2-3.
$ws.Rows("1:1").Insert()

# Populate the new top row with the double-header labels.
$ws.Range("A1").Value = "Extra header"
$ws.Range("I1").Value = "Input category"

# Move / collapse the selection like the saved workbook shows.
$ws.Range("H6").Select()
